$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-26 22:48:16"
$ws.Range("E3").Value = "2026-02-26 22:48:18"
$ws.Range("H3").Value = "'48%"
$ws.Range("E4").Value = "2026-02-26 22:48:21"
$ws.Range("H4").Value = "'81%"
$ws.Range("E5").Value = "2026-02-26 22:48:23"
$ws.Range("E6").Value = "2026-02-26 22:48:25"
$ws.Range("H6").Value = "'85%"
$ws.Range("O6").Value = "11.6 °C"
$ws.Range("E7").Value = "2026-02-26 22:48:28"
$ws.Range("E8").Value = "2026-02-26 22:48:30"
$ws.Range("E9").Value = "2026-02-26 22:48:32"
$ws.Range("O9").Value = "11.8 °C"
$ws.Range("E10").Value = "2026-02-26 22:48:33"
$ws.Range("O10").Value = "9.2 °C"
$ws.Range("E11").Value = "2026-02-26 22:48:35"
$ws.Range("O11").Value = "8.4 °C"
$ws.Range("E12").Value = "2026-02-26 22:48:36"
$ws.Range("N12").Value = "6.6 °C 22:21 TU"
$ws.Range("O12").Value = "11.1 °C"
$ws.Range("E13").Value = "2026-02-26 22:48:38"
$ws.Range("H13").Value = "'64%"
$ws.Range("O13").Value = "7.0 °C"
$ws.Range("E14").Value = "2026-02-26 22:48:41"
$ws.Range("E15").Value = "2026-02-26 22:48:43"
$ws.Range("O15").Value = "11.3 °C"
$ws.Range("E16").Value = "2026-02-26 22:48:45"
$ws.Range("O16").Value = "2.7 °C"
$ws.Range("E17").Value = "2026-02-26 22:48:48"
$ws.Range("E18").Value = "2026-02-26 22:48:50"
$ws.Range("H18").Value = "'82%"
$ws.Range("E19").Value = "2026-02-26 22:48:51"
$ws.Range("H19").Value = "'49%"
$ws.Range("E20").Value = "2026-02-26 22:48:52"
$ws.Range("E21").Value = "2026-02-26 22:48:53"
$ws.Range("J21").Value = "1027.1 hPa"
$ws.Range("O21").Value = "9.9 °C"
$ws.Range("E22").Value = "2026-02-26 22:48:56"
$ws.Range("E23").Value = "2026-02-26 22:48:58"
$ws.Range("E24").Value = "2026-02-26 22:49:01"
$ws.Range("O24").Value = "10.2 °C"
$ws.Range("E25").Value = "2026-02-26 22:49:03"
$ws.Range("K25").Value = "17.8 MJ/m2"
$ws.Range("O25").Value = "5.1 °C"
$ws.Range("E26").Value = "2026-02-26 22:49:06"
$ws.Range("E27").Value = "2026-02-26 22:49:08"
$ws.Range("O27").Value = "5.1 °C"
$ws.Range("E28").Value = "2026-02-26 22:49:11"
$ws.Range("H28").Value = "'80%"
$ws.Range("N28").Value = "4.9 °C 22:28 TU"
$ws.Range("O28").Value = "10.5 °C"
$ws.Range("E29").Value = "2026-02-26 22:49:13"
$ws.Range("N29").Value = "6.7 °C 22:27 TU"
$ws.Range("O29").Value = "11.3 °C"
$ws.Range("E30").Value = "2026-02-26 22:49:16"
$ws.Range("O30").Value = "11.9 °C"
$ws.Range("E31").Value = "2026-02-26 22:49:18"
$ws.Range("O31").Value = "11.8 °C"
$ws.Range("E32").Value = "2026-02-26 22:49:20"
$ws.Range("O32").Value = "7.6 °C"
$ws.Range("E33").Value = "2026-02-26 22:49:23"
$ws.Range("J33").Value = "1026.9 hPa"
$ws.Range("O33").Value = "8.5 °C"
$ws.Range("E34").Value = "2026-02-26 22:49:25"
$ws.Range("E35").Value = "2026-02-26 22:49:28"
$ws.Range("H35").Value = "'44%"
$ws.Range("J35").Value = "1025.6 hPa"
$ws.Range("O35").Value = "11.9 °C"
$ws.Range("E36").Value = "2026-02-26 22:49:30"
$ws.Range("N36").Value = "8.9 °C 22:26 TU"
$ws.Range("O36").Value = "12.3 °C"
$ws.Range("E37").Value = "2026-02-26 22:49:32"
$ws.Range("J37").Value = "1028.4 hPa"
$ws.Range("O37").Value = "7.5 °C"
$ws.Range("E38").Value = "2026-02-26 22:49:35"
$ws.Range("H38").Value = "'82%"
$ws.Range("O38").Value = "10.9 °C"
$ws.Range("E39").Value = "2026-02-26 22:49:37"
$ws.Range("N39").Value = "0.2 °C 22:09 TU"
$ws.Range("E40").Value = "2026-02-26 22:49:40"
$ws.Range("H40").Value = "'68%"
$ws.Range("O40").Value = "9.3 °C"
$ws.Range("E41").Value = "2026-02-26 22:49:42"
$ws.Range("O41").Value = "11.0 °C"
$ws.Range("E42").Value = "2026-02-26 22:49:44"
$ws.Range("N42").Value = "7.0 °C 22:25 TU"
$ws.Range("O42").Value = "11.1 °C"
$ws.Range("E43").Value = "2026-02-26 22:49:46"
$ws.Range("H43").Value = "'74%"
$ws.Range("E44").Value = "2026-02-26 22:49:49"
$ws.Range("E45").Value = "2026-02-26 22:49:51"
$ws.Range("J45").Value = "1026.1 hPa"
$ws.Range("O45").Value = "10.4 °C"
$ws.Range("E46").Value = "2026-02-26 22:49:53"
$ws.Range("O46").Value = "11.3 °C"
